$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Locate the picture shape named "Imagem 13" and move it up
# (its y offset changes from 2461450 EMU to 141997 EMU; x/width/height unchanged).
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Imagem 13") {
        $shp.Top = 141997 / 12700
        break
    }
}
